# edit.ps1 - reproduce the authored change:
#   1. Slide 6's table switches from the custom "Integral" table style to
#      the built-in "No Style, Table Grid" style (GUID swap in <a:tableStyleId>).
#   2. The presentation's theme colour scheme is swapped from the custom
#      "Integral" palette to the stock "Office Theme" palette (and, in the
#      OOXML, the previously-unused "Office Theme" theme part picks up the
#      "Integral" palette instead) -- i.e. the deck's Design is switched
#      from "Integral" to "Office Theme".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 6 (the table is shape 2 of that slide, but we
#    scan defensively in case shape order ever differs).
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shape = $slide6.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{494C54EB-041E-4B45-BD56-0CB8DCAAFA96}")
    }
}

# ---------------------------------------------------------------------
# 2) Switch the deck's colour scheme from "Integral" to "Office Theme".
#    ThemeColorScheme indices follow clrScheme document order:
#    1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#    8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$officeThemeColors = @{
    1  = 0        # dk1     000000
    2  = 16777215 # lt1     FFFFFF
    3  = 6968388  # dk2     44546A
    4  = 15132391 # lt2     E7E6E6
    5  = 13998939 # accent1 5B9BD5
    6  = 3243501  # accent2 ED7D31
    7  = 10855845 # accent3 A5A5A5
    8  = 49407    # accent4 FFC000
    9  = 12874308 # accent5 4472C4
    10 = 4697456  # accent6 70AD47
    11 = 12673797 # hlink   0563C1
    12 = 7491477  # folHlink 954F72
}

foreach ($idx in $officeThemeColors.Keys) {
    $colorScheme.Item($idx).RGB = $officeThemeColors[$idx]
}
